# Updated symbol list on Tue Dec 20 14:55:59 UTC 2022 with GitHub Actions
# Applies refreshed crypto price / volume-label values to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    # Cells in this sheet hold numeric-looking data as plain text.
    # Force a text number format before assignment so Excel does not
    # silently coerce the string into a floating point number, then
    # restore the cell's style so no residual formatting is left behind.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Column D (Price) updates
Set-TextValue $ws "D2"  "248.61"
Set-TextValue $ws "D3"  "22.62"
Set-TextValue $ws "D4"  "5.389"
Set-TextValue $ws "D5"  "0.05605"
Set-TextValue $ws "D7"  "6.354"
Set-TextValue $ws "D8"  "0.8143"
Set-TextValue $ws "D9"  "0.9140"
Set-TextValue $ws "D10" "0.1419"
Set-TextValue $ws "D11" "0.07494"
Set-TextValue $ws "D12" "0.03172"
Set-TextValue $ws "D13" "0.03098"
Set-TextValue $ws "D14" "0.09324"
Set-TextValue $ws "D15" "3.563"
Set-TextValue $ws "D16" "0.001634"
Set-TextValue $ws "D17" "0.04724"
Set-TextValue $ws "D18" "0.0005770"
Set-TextValue $ws "D19" "0.006395"
Set-TextValue $ws "D20" "0.004981"
Set-TextValue $ws "D21" "0.001033"
Set-TextValue $ws "D22" "0.0001503"
Set-TextValue $ws "D23" "3.729"
Set-TextValue $ws "D24" "2.168"
Set-TextValue $ws "D25" "0.3253"
Set-TextValue $ws "D26" "0.1304"
Set-TextValue $ws "D28" "0.0003005"
Set-TextValue $ws "D40" "0.03970"
Set-TextValue $ws "D41" "0.007013"
Set-TextValue $ws "D42" "0.1064"
Set-TextValue $ws "D43" "0.003406"
Set-TextValue $ws "D44" "0.007556"
Set-TextValue $ws "D45" "0.00005582"
Set-TextValue $ws "D48" "0.6010"
Set-TextValue $ws "D49" "0.2226"
Set-TextValue $ws "D50" "0.00002104"
Set-TextValue $ws "D51" "0.01012"

# Column E (Volume(1h) label) updates
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E27").Value = "26AAXTokenAABWorstin24h"
